$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.470.66"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").Value = "1.898.72"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.38%  "

# Row 5
$ws.Range("D5").Value = "238.03"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.48%  "

# Row 7
$ws.Range("D7").Value = "0.4913"
$ws.Range("E7").Value = "  +0.51%  "

# Row 8
$ws.Range("E8").Value = "  +0.65%  "

# Row 9
$ws.Range("D9").Value = "0.06717"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("D10").Value = "1.891.97"
$ws.Range("E10").Value = "  +0.02%  "

# Row 11
$ws.Range("E11").Value = "  +1.88%  "

# Row 12
$ws.Range("D12").Value = "0.07329"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13
$ws.Range("D13").Value = "5.155"
$ws.Range("E13").Value = "  +2.63%  "

# Row 14
$ws.Range("D14").Value = "87.79"
$ws.Range("E14").Value = "  -1.73%  "

# Row 15
$ws.Range("D15").Value = "0.6666"
$ws.Range("E15").Value = "  +0.49%  "

# Row 16
$ws.Range("D16").Value = "30.444.45"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17
$ws.Range("D17").Value = "13.48"
$ws.Range("E17").Value = "  +3.43%  "

# Row 18
$ws.Range("D18").Value = "0.000007853"
$ws.Range("E18").Value = "  -0.58%  "

# Row 19
$ws.Range("D19").Value = "0.9987"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("D20").Value = "2.128.43"
$ws.Range("E20").Value = "  -0.21%  "

# Row 21
$ws.Range("E21").Value = "  +12.72%  "

# Row 22
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.43%  "

# Row 23
$ws.Range("D23").Value = "191.59"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").Value = "6.129"
$ws.Range("E24").Value = "  +0.40%  "

# Row 25
$ws.Range("D25").Value = "9.495"
$ws.Range("E25").Value = "  +1.86%  "

# Row 26
$ws.Range("D26").Value = "163.42"
$ws.Range("E26").Value = "  +2.54%  "

# Row 27
$ws.Range("D27").Value = "18.27"
$ws.Range("E27").Value = "  -0.51%  "

# Row 28
$ws.Range("D28").Value = "1.938"
$ws.Range("E28").Value = "  +5.55%  "

# Row 29
$ws.Range("D29").Value = "1.458"
$ws.Range("E29").Value = "  +3.61%  "

# Row 30
$ws.Range("D30").Value = "4.350"
$ws.Range("E30").Value = "  +2.28%  "

# Row 31
$ws.Range("D31").Value = "0.09200"
$ws.Range("E31").Value = "  +2.10%  "

# Row 32
$ws.Range("D32").Value = "4.052"
$ws.Range("E32").Value = "  +2.90%  "

# Row 33
$ws.Range("D33").Value = "0.05199"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34
$ws.Range("D34").Value = "0.7417"
$ws.Range("E34").Value = "  +1.71%  "

# Row 35
$ws.Range("D35").Value = "1.105"
$ws.Range("E35").Value = "  +1.92%  "

# Row 36
$ws.Range("D36").Value = "2.706"
$ws.Range("E36").Value = "  +0.10%  "

# Row 37
$ws.Range("D37").Value = "0.01810"
$ws.Range("E37").Value = "  -0.38%  "

# Row 38
$ws.Range("D38").Value = "2.676"
$ws.Range("E38").Value = "  +0.45%  "

# Row 39
$ws.Range("D39").Value = "0.9203"
$ws.Range("E39").Value = "  -0.43%  "

# Row 40
$ws.Range("D40").Value = "2.049"
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("D41").Value = "0.4396"
$ws.Range("E41").Value = "  -0.14%  "

# Row 42
$ws.Range("D42").Value = "5.920"
$ws.Range("E42").Value = "  +3.40%  "

# Row 43
$ws.Range("D43").Value = "106.46"
$ws.Range("E43").Value = "  +1.81%  "

# Row 44
$ws.Range("D44").Value = "0.9943"
$ws.Range("E44").Value = "  -0.65%  "

# Row 45
$ws.Range("D45").Value = "69.05"
$ws.Range("E45").Value = "  +20.46%  "

# Row 46
$ws.Range("D46").Value = "0.1373"
$ws.Range("E46").Value = "  +3.03%  "

# Row 47
$ws.Range("D47").Value = "7.631"
$ws.Range("E47").Value = "  +3.92%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.001"
$ws.Range("E48").Value = "  +3.51%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "35.13"
$ws.Range("E49").Value = "  +5.53%  "

# Row 50
$ws.Range("D50").Value = "0.05829"
$ws.Range("E50").Value = "  -0.20%  "

# Row 51
$ws.Range("E51").Value = "  -3.47%  "
